# Update Name of Algo
# Apply updated imputed values to the "C" and "F"(E) columns for a few rows
# as produced by the RandomForest algorithm re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value  = -11.0661
$ws.Range("C10").Value = -12.8011
$ws.Range("C12").Value = -14.22890000000001
$ws.Range("E13").Value = 12.3504
$ws.Range("C18").Value = -14.06410000000001
